# Update "想去人数" (want-to-go count) values in the 展览 and 全部类型 sheets
# to reflect the latest scrape at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 6669
    $ws.Range("F6").Value = 2074
    $ws.Range("F7").Value = 1578
    $ws.Range("F10").Value = 461
}
